# Apply the "cryptos" price/volume refresh for the GitHub Actions update.
# Values are written as plain text (same as the source data), so any cell
# whose new text looks like a number is prefixed with a leading apostrophe
# (Excel's standard "number stored as text" quote-prefix) to stop Excel
# from silently re-typing it as a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.911.78'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').Value = '2.903.97'
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''570.23'
$ws.Range('E5').Value = '  -3.50%  '
$ws.Range('D6').Value = '''144.02'
$ws.Range('E6').Value = '  -3.10%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D9').Value = '2.902.42'
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('D10').Value = '''7.02'
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('D11').Value = '''0.146'
$ws.Range('E11').Value = '  -3.77%  '
$ws.Range('E12').Value = '  -2.13%  '
$ws.Range('D13').Value = '''0.0000230'
$ws.Range('E13').Value = '  -3.20%  '
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '3.384.93'
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').Value = '61.864.98'
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('D18').Value = '2.905.53'
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('E19').Value = '  -2.51%  '
$ws.Range('D20').Value = '''429.00'
$ws.Range('E20').Value = '  -1.91%  '
$ws.Range('D21').Value = '''12.90'
$ws.Range('E21').Value = '  -4.42%  '
$ws.Range('E22').Value = '  -2.19%  '
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('D24').Value = '''78.86'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('D25').Value = '''12.00'
$ws.Range('E25').Value = '  +0.75%  '
$ws.Range('D26').Value = '''10.15'
$ws.Range('E26').Value = '  -9.32%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '''2.02'
$ws.Range('E28').Value = '  -4.81%  '
$ws.Range('E29').Value = '  +6.79%  '
$ws.Range('E30').Value = '  -5.44%  '
$ws.Range('D31').Value = '''2.51'
$ws.Range('E31').Value = '  -3.18%  '
$ws.Range('E32').Value = '  -6.98%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  -3.54%  '
$ws.Range('D35').Value = '''25.56'
$ws.Range('E35').Value = '  -2.36%  '
$ws.Range('D36').Value = '''0.958'
$ws.Range('E36').Value = '  -2.91%  '
$ws.Range('E37').Value = '  -3.56%  '
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('D39').Value = '''2.87'
$ws.Range('E39').Value = '  -8.06%  '
$ws.Range('E40').Value = '  -6.46%  '
$ws.Range('E41').Value = '  -1.72%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').Value = '''8.14'
$ws.Range('E42').Value = '  -3.14%  '
$ws.Range('B43').Value = 'Arweave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D43').Value = '''40.91'
$ws.Range('E43').Value = '  +2.62%  '
$ws.Range('E44').Value = '  -4.40%  '
$ws.Range('D45').Value = '2.698.46'
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('E46').Value = '  -1.71%  '
$ws.Range('D47').Value = '''131.76'
$ws.Range('E47').Value = '  -3.18%  '
$ws.Range('D48').Value = '''346.50'
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').Value = '''21.53'
$ws.Range('E51').Value = '  -5.50%  '
